$d = $word.ActiveDocument

function Insert-Run {
    param($pos, $text)
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($text)
    return $pos + $text.Length
}

# Locate the insertion point: right after "ngkat keras yang digunakan " and
# before the _GoBack bookmark / the old "dalam membuat ..." run.
$rng = $d.Content
$rng.Find.Execute("ngkat keras yang digunakan ") | Out-Null
$pos = $rng.End

# Insert the new (expanded) spesifikasi text ahead of the bookmark, in the
# same run-by-run shape as the target: plain runs first, remembering the
# span that must become italic ("motion capture") so we can apply the
# formatting afterwards (avoids inheriting italics while typing).
$pos = Insert-Run $pos "dalam membuat dan melakukan uji coba adalah satu unit laptop dengan spesifikasi: CPU Intel"
$pos = Insert-Run $pos " Core"
$pos = Insert-Run $pos " I7 7700"
$pos = Insert-Run $pos " "
$pos = Insert-Run $pos "HQ, Memori 24 GB"
$pos = Insert-Run $pos " DDR4"
$pos = Insert-Run $pos ", GPU NVIDIA GTX 1060 6GB, SSD NVME SAMSUNG 120 GB, HDD SATA 1"
$pos = Insert-Run $pos " "
$pos = Insert-Run $pos "TB. Perekaman data animasi didapatkan dengan menggunakan alat "

$italicStart = $pos
$pos = Insert-Run $pos "motion capture"
$italicEnd = $pos

$pos = Insert-Run $pos " yang disediakan oleh Universitas Gunadarma. Perangkat lunak yang digunakan meliputi Git, GitHub, Mozzila Firefox, Microsoft Windows 10 Home, Visual Studio 15 2017 Community Edition, "
$pos = Insert-Run $pos "Premake 5, "

# Apply italics to "motion capture" only, after all inserts are done so the
# formatting does not leak into later runs.
$italicRange = $d.Range($italicStart, $italicEnd)
$italicRange.Font.Italic = 1

# Remove the now-duplicated old text that used to follow the _GoBack
# bookmark, leaving only the trailing "dan Microsoft Word 2016." sentence.
$searchRng = $d.Range($pos, $d.Content.End)
$searchRng.Find.Execute("dalam membuat dan melakukan uji coba adalah satu unit laptop dengan spesifikasi: CPU Intel I7 7700") | Out-Null
$oldStart = $searchRng.Start

$searchRng2 = $d.Range($pos, $d.Content.End)
$searchRng2.Find.Execute("dan Microsoft Word 2016.") | Out-Null
$oldTailStart = $searchRng2.Start

$delRange = $d.Range($oldStart, $oldTailStart)
$delRange.Delete()
